$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns store values as literal text
# (e.g. "315.81", "6.54%") rather than numbers. Apply a Text number format
# to each updated cell before writing so Excel keeps the new value as text
# instead of auto-converting it to a number/percentage.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "315.73"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "6.55%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "45.25"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "6.94%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.181"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.92%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08088"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "6.64%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.539"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.33%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.684"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "4.80%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "17.45%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "8.03%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1936"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "5.06%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09365"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "4.04%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04245"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "6.09%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1045"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.67%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001317"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "3.04%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005908"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.32%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.04%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.39%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3371"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.55%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.226"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "4.47%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1357"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-3.09%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.3147"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.91%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04279"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "5.28%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001282"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.08%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004215"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "7.60%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "9.36%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02703"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "11.83%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05460"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "4.72%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.005866"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-3.25%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007774"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.21%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1425"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "6.98%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007338"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-2.93%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008587"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "18.28%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3141"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "5.96%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006812"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.50%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.36%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05447"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "18.25%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003987"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-5.10%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002093"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.36%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001994"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.36%"
